# Apply the "Descrição de use cases revisto" edit:
# - Update the text in D9 (shared string) to remove " ou espera"
# - Update the active selection from D10 to D9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell text content
$ws.Range("D9").Value = "3. Mostra lista de carros em produção"

# Update the selected/active cell to D9 (was D10)
$ws.Range("D9").Select()
